$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$ws.Range("B9").Value = "200c320a25b192e3c83440f334527e01"
$ws.Range("B17").Value = "0ac9d1f217e8aa92141b97ddd2e464a5"
$ws.Range("B89").Value = "540c6e9b1efc86a7027d6bfbd80c73c2"
$ws.Range("B99").Value = "3ed806b97270274a88c3d0a88769021f"
$ws.Range("B110").Value = "1cbee20c6dd597308e23e402c1cb3429"
$ws.Range("B126").Value = "0e7449a6be04ef7efd69afaf0df094cc"
$ws.Range("B154").Value = "7883f0f152cc9d9bb5a1fc710f211227"
$ws.Range("B159").Value = "aaac34bf96dd1a6394dd3ca31665d0c9"
$ws.Range("B160").Value = "86c3466b53645a70143a60d23010a457"
$ws.Range("B183").Value = "0ca4f74849b3b000caf79995a5889750"
$ws.Range("B200").Value = "a84840917c81c5a306c6fab2c73dd40c"
$ws.Range("B222").Value = "60e39cb58668f837f3ef9ef35b8dd94c"
$ws.Range("B228").Value = "ac8bed145257d1de473f50ed7b7c69cc"
$ws.Range("B229").Value = "7d4adc924049e1e26eb3f440c3450a2b"
$ws.Range("B278").Value = "ff0cdaad1bb498b10fd0b974320bdfa6"
$ws.Range("B281").Value = "d47b4c2c37695aeaedf46052fc07213c"
$ws.Range("B293").Value = "ae5dcbe8cd6a13a23e310c4446ca6fc6"
$ws.Range("B335").Value = "ce0d246ac8e46bde9469712017fd6d68"
$ws.Range("B339").Value = "0cfcf0cdbc873d2da6b6d2d79315cafe"
$ws.Range("B411").Value = "6b086a7c91481ced87dd9086c965e01a"
$ws.Range("B420").Value = "0841f66eec1f7caf51680bed6f5054c6"
$ws.Range("B448").Value = "e68c149eff2080b4680818ec6449e4f1"
$ws.Range("B523").Value = "c85280c7cb5f69f7fdc4117e7b066ac0"
$ws.Range("B542").Value = "b526e2e952a95b9a09ec2a8738f95769"
$ws.Range("B561").Value = "5cbb749084cfb11e073fabbd9fa5cca4"
$ws.Range("B574").Value = "58573e77841fe155b3a08ae01bb558a7"
$ws.Range("B580").Value = "a7bcf87a3faf7a525f8737330e459fae"
$ws.Range("B592").Value = "2a0370be441331729a17ae4b1bdd77b2"
$ws.Range("B688").Value = "02796346b86ff6d9d6c7fce4bac0cac5"
$ws.Range("B693").Value = "2d3d3d86d21bacb7bbb70fb06d396780"
$ws.Range("B711").Value = "2bbbc64dc8be0d94d0befb3fe111fabd"
$ws.Range("B723").Value = "4ddd244a02ae194577a8d7a8096c1357"
$ws.Range("B764").Value = "9bf8a4ae7038adce1136408b3f7c88f8"
$ws.Range("B776").Value = "ec7cbf44da2741d451e3a0d8eb8e7bff"
$ws.Range("B819").Value = "19e459ae140fd3ca9c68c0372a062362"
$ws.Range("B823").Value = "ce02acf55c77ea096712c1a555e3035c"
$ws.Range("B824").Value = "31a7eec50e7a0a340aa5949d03d55669"
$ws.Range("B827").Value = "af8a0fdf3300e2447c7ee9846c20357a"
$ws.Range("B833").Value = "138c1287037ebf103f817fe612d3f27d"
$ws.Range("B835").Value = "820a409f29375b7c62388a0b687f0f64"
$ws.Range("B863").Value = "285c30639f51b8e86c2e6928f88b95be"
$ws.Range("B870").Value = "2868f8250a17e53d0e7b5226a008fd5f"
$ws.Range("B877").Value = "49337a52b429ecf74c50751a163f422b"
$ws.Range("B913").Value = "d94cd154838ed1d62759d2262babeefa"
